# BloomSail.pptx edit: remove the "Rechteck 15" banner/caption shape
# (the grey rounded caption box reading "Three months right in the
# middle of a cyanobacteria bloom... / ... with a sailing vessel and a
# submersible pCO2-sensor!") together with its entrance animation from
# slide 1, leaving only the research-foci pictures/connectors.

$p = $ppt.ActivePresentation
$s = $p.Slides.Item(1)

# Drop the fly-in "Appear" animation that targets the banner shape
# (spid 16) before removing the shape itself so no orphaned timing
# node is left behind referencing a shape that no longer exists.
$seq = $s.TimeLine.MainSequence
for ($i = $seq.Count; $i -ge 1; $i--) {
    $seq.Item($i).Delete()
}

# Remove the banner/caption rectangle shape entirely.
$s.Shapes.Item("Rechteck 15").Delete()
